$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: collapse the multi-line "код\n идентификатор\n PLU" label into a single line.
$ws.Range("A1").Value = "код идентификатор PLU"

# C4: rounding fix, 35 -> 8
$ws.Range("C4").Value = 8

# Column A gets an explicit width (stored width ~28 in the xlsx)
$ws.Columns.Item(1).ColumnWidth = 27.14

# Page setup: A4 portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on D14 when the file was saved
$ws.Range("D14").Select()
